$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 19:52"

# Row 4
$ws.Range("B4").Value = 300432
$ws.Range("C4").Value = 23271
$ws.Range("D4").Value = 14514
$ws.Range("E4").Value = 277764
$ws.Range("F4").Value = 7827
$ws.Range("G4").Value = 750
$ws.Range("H4").Value = 8154

# Row 8
$ws.Range("D8").Value = 15438
$ws.Range("E8").Value = 59167
$ws.Range("G8").Value = 1053
$ws.Range("H8").Value = 7560

# Row 12
$ws.Range("B12").Value = 23934
$ws.Range("C12").Value = 3013
$ws.Range("D12").Value = 786
$ws.Range("E12").Value = 22647
$ws.Range("F12").Value = 1311
$ws.Range("G12").Value = 76
$ws.Range("H12").Value = 501

# Row 63
$ws.Range("A63").Value = "Marruecos"
$ws.Range("B63").Value = 883
$ws.Range("C63").Value = 92
$ws.Range("D63").Value = 65
$ws.Range("E63").Value = 760
$ws.Range("F63").Value = 1
$ws.Range("G63").Value = 10
$ws.Range("H63").Value = 58

# Row 64
$ws.Range("A64").Value = "Irak"
$ws.Range("B64").Value = 878
$ws.Range("C64").Value = 58
$ws.Range("D64").Value = 259
$ws.Range("E64").Value = 563
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 56

# Row 65
$ws.Range("A65").Value = "Hong Kong"
$ws.Range("B65").Value = 862
$ws.Range("C65").Value = 17
$ws.Range("D65").Value = 173
$ws.Range("E65").Value = 685
$ws.Range("F65").Value = 8
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 4

# Row 91
$ws.Range("A91").Value = "Jordania"
$ws.Range("B91").Value = 323
$ws.Range("C91").Value = 13
$ws.Range("D91").Value = 74
$ws.Range("E91").Value = 244
$ws.Range("F91").Value = 5
$ws.Range("H91").Value = 5

# Row 92
$ws.Range("A92").Value = "Burkina Faso"
$ws.Range("B92").Value = 318
$ws.Range("C92").Value = 16
$ws.Range("D92").Value = 66
$ws.Range("E92").Value = 236
$ws.Range("F92").Value = 0
$ws.Range("H92").Value = 16

# Row 95
$ws.Range("E95").Value = 214
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 2

# Row 111
$ws.Range("B111").Value = 155
$ws.Range("C111").Value = 2
$ws.Range("E111").Value = 96

# Row 137
$ws.Range("B137").Value = 52
$ws.Range("C137").Value = 1
$ws.Range("E137").Value = 52

# Row 167
$ws.Range("A167").Value = "Guinea-Bisau"
$ws.Range("F167").Value = 0

# Row 168
$ws.Range("A168").Value = "Antigua y Barbuda"
$ws.Range("F168").Value = 1

# Row 173
$ws.Range("A173").Value = "Fiyi"
$ws.Range("C173").Value = 5
$ws.Range("F173").Value = 0

# Row 174
$ws.Range("A174").Value = "Granada"
$ws.Range("C174").Value = 0
$ws.Range("F174").Value = 2

# Row 175
$ws.Range("A175").Value = "Groenlandia"
$ws.Range("C175").Value = 1
$ws.Range("E175").Value = 8
$ws.Range("H175").Value = 0

# Row 176
$ws.Range("A176").Value = "Curazao"
$ws.Range("B176").Value = 11
$ws.Range("D176").Value = 3
$ws.Range("E176").Value = 7
$ws.Range("H176").Value = 1

# Row 177
$ws.Range("A177").Value = "Seychelles"

# Row 178
$ws.Range("A178").Value = "Laos"
$ws.Range("C178").Value = 0
$ws.Range("E178").Value = 10
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

# Row 181
$ws.Range("A181").Value = "Liberia"
$ws.Range("C181").Value = 3
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 9
$ws.Range("G181").Value = 1
$ws.Range("H181").Value = 1

# Row 183
$ws.Range("A183").Value = "Suazilandia"
$ws.Range("C183").Value = 0

# Row 185
$ws.Range("A185").Value = "Republica del Chad"
$ws.Range("C185").Value = 1
